# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the refreshed counts recorded in the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")

$exhibitUpdates = @{
    2  = 606
    4  = 1283
    5  = 1140
    6  = 14242
    7  = 16144
    9  = 80
    19 = 34
    24 = 6459
    25 = 968
    29 = 5672
    30 = 91
    32 = 164
    33 = 4701
    34 = 12
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# --- Sheet "全部类型" ---------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @{
    2  = 606
    4  = 1283
    5  = 1140
    6  = 14242
    7  = 16145
    9  = 80
    19 = 34
    25 = 6459
    26 = 968
    31 = 5672
    32 = 91
    34 = 164
    35 = 4701
    36 = 12
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
